# Append the new "2021年" data row (row 11) to Sheet1, mirroring the
# structure of the existing rows (row 2 = 2012年 ... row 10 = 2020年).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 11

# Column A: year label, formatted like the other year cells (bold, thin box
# border, centered/top aligned - matching the style used for A2:A10).
$ws.Cells.Item($row, 1).Value = "2021年"
$ws.Cells.Item($row, 1).Font.Bold = $true
$ws.Cells.Item($row, 1).HorizontalAlignment = -4108
$ws.Cells.Item($row, 1).VerticalAlignment = -4160
$ws.Cells.Item($row, 1).Borders.LineStyle = 1

# Column B .. AQ: the numeric indicators for 2021, column E has no reported
# value for this series (kept blank, same as it is for 2016-2020).
$ws.Cells.Item($row, 2).Value = 202.69
$ws.Cells.Item($row, 3).Value = 47.35
$ws.Cells.Item($row, 4).Value = 3.4
$ws.Cells.Item($row, 5).Value = "'"
$ws.Cells.Item($row, 5).Style = "Normal"
$ws.Cells.Item($row, 6).Value = 76.18000000000001
$ws.Cells.Item($row, 7).Value = 270.43
$ws.Cells.Item($row, 8).Value = 13.16
$ws.Cells.Item($row, 9).Value = 578.15
$ws.Cells.Item($row, 10).Value = 11.24
$ws.Cells.Item($row, 11).Value = 5061.93
$ws.Cells.Item($row, 12).Value = 2.17
$ws.Cells.Item($row, 13).Value = 6.85
$ws.Cells.Item($row, 14).Value = 0.97
$ws.Cells.Item($row, 15).Value = 11.87
$ws.Cells.Item($row, 16).Value = 64.73
$ws.Cells.Item($row, 17).Value = 4.92
$ws.Cells.Item($row, 18).Value = 3.99
$ws.Cells.Item($row, 19).Value = 37
$ws.Cells.Item($row, 20).Value = 161.02
$ws.Cells.Item($row, 21).Value = 1007.62
$ws.Cells.Item($row, 22).Value = 154.24
$ws.Cells.Item($row, 23).Value = 236.23
$ws.Cells.Item($row, 24).Value = 152.52
$ws.Cells.Item($row, 25).Value = 48.41
$ws.Cells.Item($row, 26).Value = 184.8
$ws.Cells.Item($row, 27).Value = 1.14
$ws.Cells.Item($row, 28).Value = 111.28
$ws.Cells.Item($row, 29).Value = 32.71
$ws.Cells.Item($row, 30).Value = 10.03
$ws.Cells.Item($row, 31).Value = 5.11
$ws.Cells.Item($row, 32).Value = 437.99
$ws.Cells.Item($row, 33).Value = 176.06
$ws.Cells.Item($row, 34).Value = 9.779999999999999
$ws.Cells.Item($row, 35).Value = 336.02
$ws.Cells.Item($row, 36).Value = 2.82
$ws.Cells.Item($row, 37).Value = 43.51
$ws.Cells.Item($row, 38).Value = 109.18
$ws.Cells.Item($row, 39).Value = 209.32
$ws.Cells.Item($row, 40).Value = 25.07
$ws.Cells.Item($row, 41).Value = 116.64
$ws.Cells.Item($row, 42).Value = 143.78
$ws.Cells.Item($row, 43).Value = 11.16
